# Apply the edits described by the commit "Minnor changes in OpenSource results"
#  1. Update the sheet view (scroll position / selected cell)
#  2. Switch the STDEV.S / VAR.S formulas (sample statistics) to STDEV.P / VAR.P (population statistics)
#     in two summary blocks (rows 32-33 and rows 61-62)
# Note: the workbook-level cached x15ac:absPath (the path shown in the diff header) reflects the
# real filesystem location Excel last saved the file to; it is not an addressable part of the
# documented Excel object model, so it cannot be set from script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the view: scroll to column K / row 22 and select K62 ---
$win = $excel.ActiveWindow
try {
    $win.ScrollColumn = 11
    $win.ScrollRow = 22
} catch {}

# --- 2a. Row 32 (STDEV.S -> STDEV.P) ---
$ws.Range("T32").Formula = "=STDEV.P(C29,M29)"
$ws.Range("U32").Formula = "=STDEV.P(E29,N29)"
$ws.Range("V32").Formula = "=STDEV.P(F29,O29)"
$ws.Range("W32").Formula = "=STDEV.P(H29,Q29)"
$ws.Range("X32").Formula = "=STDEV.P(L29,S29)"
$ws.Range("Y32").Formula = "=STDEV.P(J29,R29)"

# --- 2b. Row 33 (VAR.S -> VAR.P) ---
$ws.Range("T33").Formula = "=VAR.P(C29,M29)"
$ws.Range("U33").Formula = "=VAR.P(E29,N29)"
$ws.Range("V33").Formula = "=VAR.P(F29,O29)"
$ws.Range("W33").Formula = "=VAR.P(H29,Q29)"
$ws.Range("X33").Formula = "=VAR.P(L29,S29)"
$ws.Range("Y33").Formula = "=VAR.P(J29,R29)"

# --- 2c. Row 61 (STDEV.S -> STDEV.P, formulas un-shared) ---
$ws.Range("K61").Formula = "=STDEV.P(C58,G58)"
$ws.Range("L61").Formula = "=STDEV.P(D58,H58)"
$ws.Range("M61").Formula = "=STDEV.P(E58,I58)"
$ws.Range("N61").Formula = "=STDEV.P(F58,J58)"

# --- 2d. Row 62 (VAR.S -> VAR.P, formulas un-shared) ---
$ws.Range("K62").Formula = "=VAR.P(C58,G58)"
$ws.Range("L62").Formula = "=VAR.P(D58,H58)"
$ws.Range("M62").Formula = "=VAR.P(E58,I58)"
$ws.Range("N62").Formula = "=VAR.P(F58,J58)"

# --- 1 (cont.) select the new active cell, matching the recorded selection in the diff ---
$ws.Range("K62").Select()
